# Auto-generated Excel COM-interop script to update Leve profit values
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 229.2
$ws.Range("I2").Value = 229.2
$ws.Range("K2").Value = 229.2
$ws.Range("M2").Value = -116.2

$ws.Range("H4").Value = 649.2
$ws.Range("I4").Value = 239.5
$ws.Range("K4").Value = 239.5
$ws.Range("M4").Value = -125.5

$ws.Range("H19").Value = 1527
$ws.Range("I19").Value = 955.625
$ws.Range("J19").Value = 2098.375
$ws.Range("K19").Value = 955.625
$ws.Range("L19").Value = 2098.375
$ws.Range("M19").Value = -780.625
$ws.Range("N19").Value = -2448.375

$ws.Range("H29").Value = 4998.9673
$ws.Range("I29").Value = 4998
$ws.Range("J29").Value = 4999
$ws.Range("K29").Value = 14994
$ws.Range("L29").Value = 14997
$ws.Range("M29").Value = -14713
$ws.Range("N29").Value = -15559

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 15750
$ws.Range("I6").Value = 26500
$ws.Range("K6").Value = 26500
$ws.Range("M6").Value = -26327

$ws.Range("H32").Value = 1695.4021
$ws.Range("I32").Value = 1313.2222
$ws.Range("K32").Value = 1313.2222
$ws.Range("M32").Value = -1026.2222

$ws.Range("H34").Value = 43749.5
$ws.Range("I34").Value = 44999.332
$ws.Range("J34").Value = 40000
$ws.Range("K34").Value = 44999.332
$ws.Range("L34").Value = 40000
$ws.Range("M34").Value = -44728.332
$ws.Range("N34").Value = -40542

$ws.Range("H44").Value = 49830.6
$ws.Range("J44").Value = 81669.664
$ws.Range("L44").Value = 81669.664
$ws.Range("N44").Value = -82645.664

$ws.Range("H74").Value = 1561.0233
$ws.Range("I74").Value = 1638.95
$ws.Range("K74").Value = 1638.95
$ws.Range("M74").Value = -764.95

$ws.Range("H77").Value = 1561.0233
$ws.Range("I77").Value = 1638.95
$ws.Range("K77").Value = 8194.75
$ws.Range("M77").Value = -3826.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 710.125
$ws.Range("J20").Value = 720.3
$ws.Range("L20").Value = 720.3
$ws.Range("N20").Value = -1214.3

$ws.Range("H80").Value = 793.5
$ws.Range("J80").Value = 1040.5
$ws.Range("L80").Value = 1040.5
$ws.Range("N80").Value = -3036.5

$ws.Range("H83").Value = 793.5
$ws.Range("J83").Value = 1040.5
$ws.Range("L83").Value = 5202.5
$ws.Range("N83").Value = -15186.5

$ws.Range("H134").Value = 9563.596
$ws.Range("I134").Value = 9295.656999999999
$ws.Range("K134").Value = 27886.971
$ws.Range("M134").Value = -25351.971

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2413.1667
$ws.Range("I4").Value = 896
$ws.Range("J4").Value = 9999
$ws.Range("K4").Value = 896
$ws.Range("L4").Value = 9999
$ws.Range("M4").Value = -784
$ws.Range("N4").Value = -10223

$ws.Range("H15").Value = 19999
$ws.Range("J15").Value = 19999
$ws.Range("L15").Value = 19999
$ws.Range("N15").Value = -20339

$ws.Range("H58").Value = 3336.077
$ws.Range("I58").Value = 3119.889
$ws.Range("K58").Value = 3119.889
$ws.Range("M58").Value = -2916.889

$ws.Range("H136").Value = 3336.077
$ws.Range("I136").Value = 3119.889
$ws.Range("K136").Value = 9359.667000000001
$ws.Range("M136").Value = -6809.667000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 9950
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 9950
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 29850
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -30426

$ws.Range("H48").Value = 5133
$ws.Range("J48").Value = 6199.5
$ws.Range("L48").Value = 18598.5
$ws.Range("N48").Value = -19098.5

$ws.Range("H54").Value = 8745
$ws.Range("J54").Value = 7993.6665
$ws.Range("L54").Value = 23980.9995
$ws.Range("N54").Value = -25098.9995

$ws.Range("H109").Value = 1434.2667
$ws.Range("I109").Value = 1231.8462
$ws.Range("K109").Value = 3695.5386
$ws.Range("M109").Value = -2655.5386

$ws.Range("H123").Value = 85
$ws.Range("I123").Value = 85
$ws.Range("K123").Value = 255
$ws.Range("M123").Value = 2195

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 119.10345
$ws.Range("I2").Value = 130.23529
$ws.Range("J2").Value = 103.333336
$ws.Range("K2").Value = 130.23529
$ws.Range("L2").Value = 103.333336
$ws.Range("M2").Value = -17.23528999999999
$ws.Range("N2").Value = -329.333336

$ws.Range("H15").Value = 42666.668
$ws.Range("I15").Value = 50000
$ws.Range("J15").Value = 39000
$ws.Range("K15").Value = 50000
$ws.Range("L15").Value = 39000
$ws.Range("M15").Value = -49712
$ws.Range("N15").Value = -39576

$ws.Range("H43").Value = 16086
$ws.Range("I43").Value = 999.5
$ws.Range("J43").Value = 31172.5
$ws.Range("K43").Value = 999.5
$ws.Range("L43").Value = 31172.5
$ws.Range("M43").Value = -848.5
$ws.Range("N43").Value = -31474.5

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H80").Value = 1848.75
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H81").Value = 42666.668
$ws.Range("I81").Value = 50000
$ws.Range("J81").Value = 39000
$ws.Range("K81").Value = 50000
$ws.Range("L81").Value = 39000
$ws.Range("M81").Value = -49002
$ws.Range("N81").Value = -40996

$ws.Range("H83").Value = 1848.75
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H84").Value = 42666.668
$ws.Range("I84").Value = 50000
$ws.Range("J84").Value = 39000
$ws.Range("K84").Value = 150000
$ws.Range("L84").Value = 117000
$ws.Range("M84").Value = -145008
$ws.Range("N84").Value = -126984

$ws.Range("H97").Value = 824.5
$ws.Range("I97").Value = 800
$ws.Range("J97").Value = 836.75
$ws.Range("K97").Value = 800
$ws.Range("L97").Value = 836.75
$ws.Range("M97").Value = -304
$ws.Range("N97").Value = -1828.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 3500
$ws.Range("I20").Value = 3500
$ws.Range("K20").Value = 3500
$ws.Range("M20").Value = -3274

$ws.Range("H22").Value = 2967.4443
$ws.Range("J22").Value = 2967.4443
$ws.Range("L22").Value = 2967.4443
$ws.Range("N22").Value = -3557.4443

$ws.Range("H27").Value = 2967.4443
$ws.Range("J27").Value = 2967.4443
$ws.Range("L27").Value = 2967.4443
$ws.Range("N27").Value = -3181.4443

$ws.Range("H136").Value = 2805.5625
$ws.Range("I136").Value = 1668.8889
$ws.Range("K136").Value = 5006.6667
$ws.Range("M136").Value = -2456.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3935.5715
$ws.Range("J81").Value = 6651.3335
$ws.Range("L81").Value = 13302.667
$ws.Range("N81").Value = -15424.667

$ws.Range("H84").Value = 3935.5715
$ws.Range("J84").Value = 6651.3335
$ws.Range("L84").Value = 66513.33499999999
$ws.Range("N84").Value = -77121.33499999999
